$d = $word.ActiveDocument

# Locate the heading paragraph "Dokumentation des Betriebs" (the actual
# section heading, not its occurrence inside the table of contents) and
# remove the two paragraphs that directly follow it:
#   - "Die grafische Oberfläche wir mit Hilfe des CSS-Frameworks Bootstrap ..."
#   - "Usernamen. Daraufhin werden die letzten drei Einträge ..."

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Dokumentation des Betriebs") {
        $target = $p
    }
}

if ($target -ne $null) {
    $startPara = $target.Next()
    $endPara = $startPara.Next()

    $start = $startPara.Range.Start
    $end = $endPara.Range.End

    $r = $d.Range($start, $end)
    $r.Delete()
}
